$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1404.2046
$ws.Range("I15").Value = 1404.2046
$ws.Range("K15").Value = 4212.6138
$ws.Range("M15").Value = -4043.6138
$ws.Range("H19").Value = 5058.1
$ws.Range("J19").Value = 7248.6665
$ws.Range("L19").Value = 7248.6665
$ws.Range("N19").Value = -7598.6665
$ws.Range("H33").Value = 10603.667
$ws.Range("I33").Value = 14504.866
$ws.Range("J33").Value = 850.6667
$ws.Range("K33").Value = 14504.866
$ws.Range("L33").Value = 850.6667
$ws.Range("M33").Value = -14275.866
$ws.Range("N33").Value = -1308.6667
$ws.Range("H43").Value = 228262.11
$ws.Range("J43").Value = 405297
$ws.Range("L43").Value = 405297
$ws.Range("N43").Value = -405435
$ws.Range("H51").Value = 4066
$ws.Range("I51").Value = 3200
$ws.Range("K51").Value = 3200
$ws.Range("M51").Value = -2716
$ws.Range("H69").Value = 5500
$ws.Range("J69").Value = 6000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 5500
$ws.Range("J72").Value = 6000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62736
$ws.Range("H116").Value = 8443.182000000001
$ws.Range("I116").Value = 7148
$ws.Range("K116").Value = 7148
$ws.Range("M116").Value = -3706
$ws.Range("H125").Value = 9439.286
$ws.Range("I125").Value = 10994
$ws.Range("K125").Value = 98946
$ws.Range("M125").Value = -96486
$ws.Range("H132").Value = 3192.6042
$ws.Range("I132").Value = 3519.325
$ws.Range("K132").Value = 10557.975
$ws.Range("M132").Value = -8027.974999999999
$ws.Range("H135").Value = 585.55554
$ws.Range("I135").Value = 558.75
$ws.Range("K135").Value = 5028.75
$ws.Range("M135").Value = -2493.75
$ws.Range("H138").Value = 29414450
$ws.Range("I138").Value = 1333.4615
$ws.Range("J138").Value = 47622570
$ws.Range("K138").Value = 4000.3845
$ws.Range("L138").Value = 142867710
$ws.Range("M138").Value = 1139.6155
$ws.Range("N138").Value = -142877990

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3375.25
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H132").Value = 23437.416
$ws.Range("I132").Value = 1803.5964
$ws.Range("K132").Value = 5410.789199999999
$ws.Range("M132").Value = -2880.789199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H86").Value = 9925.559999999999
$ws.Range("I86").Value = 1666.5555
$ws.Range("J86").Value = 31163
$ws.Range("K86").Value = 1666.5555
$ws.Range("L86").Value = 31163
$ws.Range("M86").Value = -543.5554999999999
$ws.Range("N86").Value = -33409
$ws.Range("H89").Value = 9925.559999999999
$ws.Range("I89").Value = 1666.5555
$ws.Range("J89").Value = 31163
$ws.Range("K89").Value = 8332.7775
$ws.Range("L89").Value = 155815
$ws.Range("M89").Value = -2716.7775
$ws.Range("N89").Value = -167047
$ws.Range("H105").Value = 20671.334
$ws.Range("I105").Value = 22805.6
$ws.Range("K105").Value = 22805.6
$ws.Range("M105").Value = -21058.6
$ws.Range("H107").Value = 2830.037
$ws.Range("I107").Value = 2597.68
$ws.Range("J107").Value = 5734.5
$ws.Range("K107").Value = 2597.68
$ws.Range("L107").Value = 5734.5
$ws.Range("M107").Value = -677.6799999999998
$ws.Range("N107").Value = -9574.5
$ws.Range("H134").Value = 852.6
$ws.Range("I134").Value = 852.6
$ws.Range("K134").Value = 2557.8
$ws.Range("M134").Value = -22.80000000000018
$ws.Range("H140").Value = 89332.664
$ws.Range("J140").Value = 89332.664
$ws.Range("L140").Value = 89332.664
$ws.Range("N140").Value = -99692.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1861.4706
$ws.Range("I58").Value = 1037.909
$ws.Range("J58").Value = 3371.3333
$ws.Range("K58").Value = 1037.909
$ws.Range("L58").Value = 3371.3333
$ws.Range("M58").Value = -834.9090000000001
$ws.Range("N58").Value = -3777.3333
$ws.Range("H74").Value = 252500
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -501748
$ws.Range("H77").Value = 252500
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 1500000
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -1508736
$ws.Range("H122").Value = 3549.2856
$ws.Range("I122").Value = 4400
$ws.Range("K122").Value = 13200
$ws.Range("M122").Value = -10750
$ws.Range("H132").Value = 5381.6113
$ws.Range("I132").Value = 5294.2666
$ws.Range("J132").Value = 5818.3335
$ws.Range("K132").Value = 15882.7998
$ws.Range("L132").Value = 17455.0005
$ws.Range("M132").Value = -13352.7998
$ws.Range("N132").Value = -22515.0005
$ws.Range("H136").Value = 1861.4706
$ws.Range("I136").Value = 1037.909
$ws.Range("J136").Value = 3371.3333
$ws.Range("K136").Value = 3113.727
$ws.Range("L136").Value = 10113.9999
$ws.Range("M136").Value = -563.7270000000003
$ws.Range("N136").Value = -15213.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 60000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H70").Value = 6141.9165
$ws.Range("I70").Value = 5626.857
$ws.Range("K70").Value = 5626.857
$ws.Range("M70").Value = -5356.857
$ws.Range("H73").Value = 6141.9165
$ws.Range("I73").Value = 5626.857
$ws.Range("K73").Value = 5626.857
$ws.Range("M73").Value = -4690.857
$ws.Range("H93").Value = 38597.8
$ws.Range("J93").Value = 41997.25
$ws.Range("L93").Value = 41997.25
$ws.Range("N93").Value = -45741.25
$ws.Range("H95").Value = 413413
$ws.Range("J95").Value = 413413
$ws.Range("L95").Value = 413413
$ws.Range("N95").Value = -418905
$ws.Range("H107").Value = 624.5454999999999
$ws.Range("I107").Value = 563.3333
$ws.Range("K107").Value = 563.3333
$ws.Range("M107").Value = 1356.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3957
$ws.Range("I7").Value = 3150.3684
$ws.Range("K7").Value = 3150.3684
$ws.Range("M7").Value = -3038.3684
$ws.Range("H40").Value = 4040.1538
$ws.Range("I40").Value = 3370.8
$ws.Range("J40").Value = 4458.5
$ws.Range("K40").Value = 3370.8
$ws.Range("L40").Value = 4458.5
$ws.Range("M40").Value = -3234.8
$ws.Range("N40").Value = -4730.5
$ws.Range("H61").Value = 3057.8948
$ws.Range("I61").Value = 2783.3333
$ws.Range("K61").Value = 2783.3333
$ws.Range("M61").Value = -2581.3333
$ws.Range("H107").Value = 10959.25
$ws.Range("I107").Value = 10959.25
$ws.Range("K107").Value = 10959.25
$ws.Range("M107").Value = -9039.25
$ws.Range("H113").Value = 3057.8948
$ws.Range("I113").Value = 2783.3333
$ws.Range("K113").Value = 2783.3333
$ws.Range("M113").Value = -613.3332999999998
$ws.Range("H122").Value = 3663.25
$ws.Range("I122").Value = 3379.125
$ws.Range("J122").Value = 4799.75
$ws.Range("K122").Value = 10137.375
$ws.Range("L122").Value = 14399.25
$ws.Range("M122").Value = -7687.375
$ws.Range("N122").Value = -19299.25
$ws.Range("H126").Value = 3957
$ws.Range("I126").Value = 3150.3684
$ws.Range("K126").Value = 9451.1052
$ws.Range("M126").Value = -6981.1052
$ws.Range("H133").Value = 87775.5
$ws.Range("J133").Value = 87775.5
$ws.Range("L133").Value = 87775.5
$ws.Range("N133").Value = -92835.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 500009900
$ws.Range("J41").Value = 19800
$ws.Range("L41").Value = 19800
$ws.Range("N41").Value = -20580
$ws.Range("H62").Value = 8342.714
$ws.Range("I62").Value = 6499
$ws.Range("K62").Value = 6499
$ws.Range("M62").Value = -5875
$ws.Range("H65").Value = 8342.714
$ws.Range("I65").Value = 6499
$ws.Range("K65").Value = 32495
$ws.Range("M65").Value = -29375
$ws.Range("H81").Value = 3354.8572
$ws.Range("I81").Value = 3112.923
$ws.Range("J81").Value = 6500
$ws.Range("K81").Value = 6225.846
$ws.Range("L81").Value = 13000
$ws.Range("M81").Value = -5164.846
$ws.Range("N81").Value = -15122
$ws.Range("H84").Value = 3354.8572
$ws.Range("I84").Value = 3112.923
$ws.Range("J84").Value = 6500
$ws.Range("K84").Value = 31129.23
$ws.Range("L84").Value = 65000
$ws.Range("M84").Value = -25825.23
$ws.Range("N84").Value = -75608
$ws.Range("H107").Value = 1315.9286
$ws.Range("I107").Value = 1147.5454
$ws.Range("J107").Value = 1933.3334
$ws.Range("K107").Value = 3442.6362
$ws.Range("L107").Value = 5800.0002
$ws.Range("M107").Value = -1522.6362
$ws.Range("N107").Value = -9640.0002
$ws.Range("H113").Value = 1517.7142
$ws.Range("J113").Value = 1365
$ws.Range("L113").Value = 4095
$ws.Range("N113").Value = -8435
$ws.Range("H126").Value = 4044.3333
$ws.Range("J126").Value = 4149.9
$ws.Range("L126").Value = 12449.7
$ws.Range("N126").Value = -17389.7
